$d = $word.ActiveDocument

# The bullet currently reads:
#   "Proficiency in Information Theory and/or Network Science (e.g., graph theory)."
# It should read:
#   "Proficiency in Information Theory, Network Science (e.g., graph theory), and/or Control Theory."
#
# Replace the tail of the sentence (everything after "Information Theory") in one
# targeted Find/Replace so the full stop and new "Control Theory" clause land correctly.
$range = $d.Content
$found = $range.Find.Execute( `
    " and/or Network Science (e.g., graph theory).", `
    $true, `
    $false, `
    $false, `
    $false, `
    $false, `
    $true, `
    1, `
    $false, `
    ", Network Science (e.g., graph theory), and/or Control Theory.", `
    2)

if (-not $found) {
    throw "Could not find the 'Information Theory and/or Network Science' sentence to update."
}
